$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swaps/cycles among F:V (match detail columns); A-E (Indice..data_partida) stay fixed per row ---
$ws.Range("F125").Value = "Brescia"
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = "Cremonese"
$ws.Range("I125").Value = 3
$ws.Range("J125").Value = 3.29
$ws.Range("K125").Value = "08/11/2023 18:42"
$ws.Range("L125").Value = 4.27
$ws.Range("M125").Value = "12/11/2023 16:14"
$ws.Range("N125").Value = 3.36
$ws.Range("O125").Value = "08/11/2023 18:42"
$ws.Range("P125").Value = 3.33
$ws.Range("Q125").Value = "12/11/2023 16:14"
$ws.Range("R125").Value = 2.29
$ws.Range("S125").Value = "08/11/2023 18:42"
$ws.Range("T125").Value = 2.02
$ws.Range("U125").Value = "12/11/2023 16:14"
$ws.Range("V125").Value = "https://www.betexplorer.com/football/italy/serie-b/brescia-cremonese/xOSuO07i/"
$ws.Range("F126").Value = "Lecco"
$ws.Range("G126").Value = 3
$ws.Range("H126").Value = "Parma"
$ws.Range("I126").Value = 2
$ws.Range("J126").Value = 4.45
$ws.Range("K126").Value = "08/11/2023 18:42"
$ws.Range("L126").Value = 4.33
$ws.Range("M126").Value = "12/11/2023 16:11"
$ws.Range("N126").Value = 3.7
$ws.Range("O126").Value = "08/11/2023 18:42"
$ws.Range("P126").Value = 3.48
$ws.Range("Q126").Value = "12/11/2023 16:11"
$ws.Range("R126").Value = 1.85
$ws.Range("S126").Value = "08/11/2023 18:42"
$ws.Range("T126").Value = 1.96
$ws.Range("U126").Value = "12/11/2023 16:11"
$ws.Range("V126").Value = "https://www.betexplorer.com/football/italy/serie-b/lecco-parma/OzLlpIz9/"
$ws.Range("F127").Value = "Palermo"
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = "Cittadella"
$ws.Range("I127").Value = 1
$ws.Range("J127").Value = 1.76
$ws.Range("K127").Value = "08/11/2023 18:42"
$ws.Range("L127").Value = 1.89
$ws.Range("M127").Value = "12/11/2023 16:08"
$ws.Range("N127").Value = 3.78
$ws.Range("O127").Value = "08/11/2023 18:42"
$ws.Range("P127").Value = 3.49
$ws.Range("Q127").Value = "12/11/2023 16:08"
$ws.Range("R127").Value = 4.96
$ws.Range("S127").Value = "08/11/2023 18:42"
$ws.Range("T127").Value = 4.66
$ws.Range("U127").Value = "12/11/2023 16:12"
$ws.Range("V127").Value = "https://www.betexplorer.com/football/italy/serie-b/palermo-cittadella/W0Odrd5L/"
$ws.Range("F128").Value = "Spezia"
$ws.Range("G128").Value = 2
$ws.Range("H128").Value = "Ternana"
$ws.Range("I128").Value = 2
$ws.Range("J128").Value = 1.88
$ws.Range("K128").Value = "08/11/2023 18:42"
$ws.Range("L128").Value = 2.03
$ws.Range("M128").Value = "12/11/2023 16:14"
$ws.Range("N128").Value = 3.7
$ws.Range("O128").Value = "08/11/2023 18:42"
$ws.Range("P128").Value = 3.25
$ws.Range("Q128").Value = "12/11/2023 16:01"
$ws.Range("R128").Value = 4.29
$ws.Range("S128").Value = "08/11/2023 18:42"
$ws.Range("T128").Value = 4.39
$ws.Range("U128").Value = "12/11/2023 16:14"
$ws.Range("V128").Value = "https://www.betexplorer.com/football/italy/serie-b/spezia-ternana/EPD1sGKR/"
$ws.Range("F131").Value = "Bari"
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = "Venezia"
$ws.Range("I131").Value = 3
$ws.Range("J131").Value = 2.43
$ws.Range("K131").Value = "12/11/2023 15:42"
$ws.Range("L131").Value = 2.76
$ws.Range("M131").Value = "25/11/2023 13:52"
$ws.Range("N131").Value = 3.28
$ws.Range("O131").Value = "12/11/2023 15:42"
$ws.Range("P131").Value = 2.95
$ws.Range("Q131").Value = "25/11/2023 13:25"
$ws.Range("R131").Value = 3.1
$ws.Range("S131").Value = "12/11/2023 15:42"
$ws.Range("T131").Value = 3.06
$ws.Range("U131").Value = "25/11/2023 13:52"
$ws.Range("V131").Value = "https://www.betexplorer.com/football/italy/serie-b/bari-venezia/Ob3yyEce/"
$ws.Range("F132").Value = "Parma"
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = "Modena"
$ws.Range("I132").Value = 1
$ws.Range("J132").Value = 1.92
$ws.Range("K132").Value = "12/11/2023 16:43"
$ws.Range("L132").Value = 2
$ws.Range("M132").Value = "25/11/2023 13:56"
$ws.Range("N132").Value = 3.49
$ws.Range("O132").Value = "12/11/2023 16:43"
$ws.Range("P132").Value = 3.51
$ws.Range("Q132").Value = "25/11/2023 13:59"
$ws.Range("R132").Value = 4.29
$ws.Range("S132").Value = "12/11/2023 16:43"
$ws.Range("T132").Value = 4.12
$ws.Range("U132").Value = "25/11/2023 13:56"
$ws.Range("V132").Value = "https://www.betexplorer.com/football/italy/serie-b/parma-modena/EwgzdFZK/"
$ws.Range("F133").Value = "Cittadella"
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = "Sudtirol"
$ws.Range("I133").Value = 1
$ws.Range("J133").Value = 2.03
$ws.Range("K133").Value = "12/11/2023 16:43"
$ws.Range("L133").Value = 2.25
$ws.Range("M133").Value = "25/11/2023 13:56"
$ws.Range("N133").Value = 3.23
$ws.Range("O133").Value = "12/11/2023 16:43"
$ws.Range("P133").Value = 2.83
$ws.Range("Q133").Value = "25/11/2023 13:56"
$ws.Range("R133").Value = 4.2
$ws.Range("S133").Value = "12/11/2023 16:43"
$ws.Range("T133").Value = 4.35
$ws.Range("U133").Value = "25/11/2023 13:56"
$ws.Range("V133").Value = "https://www.betexplorer.com/football/italy/serie-b/cittadella-sudtirol/48BqZiS7/"
$ws.Range("F134").Value = "Cremonese"
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = "Lecco"
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 1.31
$ws.Range("K134").Value = "12/11/2023 16:43"
$ws.Range("L134").Value = 1.42
$ws.Range("M134").Value = "25/11/2023 13:53"
$ws.Range("N134").Value = 5.52
$ws.Range("O134").Value = "12/11/2023 16:43"
$ws.Range("P134").Value = 5.04
$ws.Range("Q134").Value = "25/11/2023 13:53"
$ws.Range("R134").Value = 10.16
$ws.Range("S134").Value = "12/11/2023 16:43"
$ws.Range("T134").Value = 7.75
$ws.Range("U134").Value = "25/11/2023 13:53"
$ws.Range("V134").Value = "https://www.betexplorer.com/football/italy/serie-b/cremonese-lecco/d2bVdeKE/"
$ws.Range("F135").Value = "Como"
$ws.Range("G135").Value = 2
$ws.Range("H135").Value = "FeralpiSalo"
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = 1.58
$ws.Range("K135").Value = "12/11/2023 15:42"
$ws.Range("L135").Value = 1.81
$ws.Range("M135").Value = "25/11/2023 13:58"
$ws.Range("N135").Value = 4.13
$ws.Range("O135").Value = "12/11/2023 15:42"
$ws.Range("P135").Value = 3.52
$ws.Range("Q135").Value = "25/11/2023 13:59"
$ws.Range("R135").Value = 6.22
$ws.Range("S135").Value = "12/11/2023 15:42"
$ws.Range("T135").Value = 5.2
$ws.Range("U135").Value = "25/11/2023 13:59"
$ws.Range("V135").Value = "https://www.betexplorer.com/football/italy/serie-b/como-feralpisalo/IBaRcy58/"
$ws.Range("F137").Value = "Catanzaro"
$ws.Range("G137").Value = 2
$ws.Range("H137").Value = "Cosenza"
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1.88
$ws.Range("K137").Value = "12/11/2023 15:42"
$ws.Range("L137").Value = 2.29
$ws.Range("M137").Value = "26/11/2023 16:14"
$ws.Range("N137").Value = 3.65
$ws.Range("O137").Value = "12/11/2023 15:42"
$ws.Range("P137").Value = 3.18
$ws.Range("Q137").Value = "26/11/2023 16:14"
$ws.Range("R137").Value = 4.37
$ws.Range("S137").Value = "12/11/2023 15:42"
$ws.Range("T137").Value = 3.63
$ws.Range("U137").Value = "26/11/2023 16:14"
$ws.Range("V137").Value = "https://www.betexplorer.com/football/italy/serie-b/catanzaro-cosenza/6NhtzYC1/"
$ws.Range("F138").Value = "Ternana"
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = "Palermo"
$ws.Range("I138").Value = 1
$ws.Range("J138").Value = 2.91
$ws.Range("K138").Value = "12/11/2023 16:43"
$ws.Range("L138").Value = 3.26
$ws.Range("M138").Value = "26/11/2023 16:14"
$ws.Range("N138").Value = 3.39
$ws.Range("O138").Value = "12/11/2023 16:43"
$ws.Range("P138").Value = 3.26
$ws.Range("Q138").Value = "26/11/2023 16:12"
$ws.Range("R138").Value = 2.5
$ws.Range("S138").Value = "12/11/2023 16:43"
$ws.Range("T138").Value = 2.42
$ws.Range("U138").Value = "26/11/2023 16:13"
$ws.Range("V138").Value = "https://www.betexplorer.com/football/italy/serie-b/ternana-palermo/OCrAkDCl/"
$ws.Range("F142").Value = "Venezia"
$ws.Range("G142").Value = 3
$ws.Range("H142").Value = "Ascoli"
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = 1.73
$ws.Range("K142").Value = "25/11/2023 14:12"
$ws.Range("L142").Value = 1.75
$ws.Range("M142").Value = "02/12/2023 13:59"
$ws.Range("N142").Value = 3.84
$ws.Range("O142").Value = "25/11/2023 14:12"
$ws.Range("P142").Value = 3.54
$ws.Range("Q142").Value = "02/12/2023 13:59"
$ws.Range("R142").Value = 5.1
$ws.Range("S142").Value = "25/11/2023 14:12"
$ws.Range("T142").Value = 5.71
$ws.Range("U142").Value = "02/12/2023 13:59"
$ws.Range("V142").Value = "https://www.betexplorer.com/football/italy/serie-b/venezia-ascoli/zFUkjp3a/"
$ws.Range("F146").Value = "Pisa"
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = "Cremonese"
$ws.Range("I146").Value = 0
$ws.Range("J146").Value = 3.16
$ws.Range("K146").Value = "25/11/2023 18:13"
$ws.Range("L146").Value = 3.52
$ws.Range("M146").Value = "02/12/2023 13:55"
$ws.Range("N146").Value = 3.37
$ws.Range("O146").Value = "25/11/2023 18:13"
$ws.Range("P146").Value = 3.33
$ws.Range("Q146").Value = "02/12/2023 13:55"
$ws.Range("R146").Value = 2.35
$ws.Range("S146").Value = "25/11/2023 18:13"
$ws.Range("T146").Value = 2.25
$ws.Range("U146").Value = "02/12/2023 13:55"
$ws.Range("V146").Value = "https://www.betexplorer.com/football/italy/serie-b/pisa-cremonese/pQywgnJt/"

# --- New row 147 (brand new match added) ---
$ws.Range("A147").Value = 146
$ws.Range("B147").Value = "italy"
$ws.Range("C147").Value = "serie-b"
$ws.Range("D147").Value = "2023-2024"
$ws.Range("E147").Value = 45262.67708333334
$ws.Range("F147").Value = "Spezia"
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = "Parma"
$ws.Range("I147").Value = 1
$ws.Range("J147").Value = 2.38
$ws.Range("K147").Value = "25/11/2023 14:12"
$ws.Range("L147").Value = 2.57
$ws.Range("M147").Value = "02/12/2023 16:14"
$ws.Range("N147").Value = 3.33
$ws.Range("O147").Value = "25/11/2023 14:12"
$ws.Range("P147").Value = 3.23
$ws.Range("Q147").Value = "02/12/2023 16:09"
$ws.Range("R147").Value = 3.14
$ws.Range("S147").Value = "25/11/2023 14:12"
$ws.Range("T147").Value = 3.04
$ws.Range("U147").Value = "02/12/2023 16:14"
$ws.Range("V147").Value = "https://www.betexplorer.com/football/italy/serie-b/spezia-parma/vmysh6Yn/"

# Copy formatting (cell styles) from row 146 into new row 147 so A147/E147 get correct number styles
$ws.Range("A146:V146").Copy() | Out-Null
$ws.Range("A147:V147").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
